$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.08817155729215642
$ws.Range("D2").Value = 0.02271376096607725
$ws.Range("E2").Value = 0.1441052394679616
$ws.Range("F2").Value = 0.6313750462135914
$ws.Range("G2").Value = 0.002411527378197862
$ws.Range("I2").Value = 0.4759696514830423
$ws.Range("K2").Value = 0.588029134313075
$ws.Range("M2").Value = 0.273934438043284
$ws.Range("O2").Value = 2.120597587544637
$ws.Range("B3").Value = 0.07770287929778874
$ws.Range("D3").Value = 0.02101260603748045
$ws.Range("E3").Value = 0.1367952155266465
$ws.Range("F3").Value = 0.6280256792361527
$ws.Range("G3").Value = 0.002414205736781145
$ws.Range("I3").Value = 0.4840765391992443
$ws.Range("K3").Value = 0.51321165650495
$ws.Range("M3").Value = 0.2436039007042794
$ws.Range("O3").Value = 2.123448152316826
$ws.Range("B4").Value = 0.07126632103182828
$ws.Range("D4").Value = 0.01996057790312733
$ws.Range("E4").Value = 0.1324338438618895
$ws.Range("F4").Value = 0.6264482697331246
$ws.Range("G4").Value = 0.002415937941384918
$ws.Range("I4").Value = 0.4893424908799058
$ws.Range("K4").Value = 0.4670949846998553
$ws.Range("M4").Value = 0.2250273761360688
$ws.Range("O4").Value = 2.126785601652074
$ws.Range("B5").Value = 0.06864137448025076
$ws.Range("D5").Value = 0.01953000646658865
$ws.Range("E5").Value = 0.1306882255653434
$ws.Range("F5").Value = 0.6259257629249007
$ws.Range("G5").Value = 0.00241666594385587
$ws.Range("I5").Value = 0.4915608597359049
$ws.Range("K5").Value = 0.4482581247519306
$ws.Range("M5").Value = 0.2174689952313003
$ws.Range("O5").Value = 2.128544068496581
$ws.Range("B6").Value = 0.06820538949649801
$ws.Range("D6").Value = 0.01945839879530098
$ws.Range("E6").Value = 0.1304002724706308
$ws.Range("F6").Value = 0.6258462617318372
$ws.Range("G6").Value = 0.002416788165727447
$ws.Range("I6").Value = 0.4919335928151911
$ws.Range("K6").Value = 0.4451276520024692
$ws.Range("M6").Value = 0.2162146387627146
$ws.Range("O6").Value = 2.128860104729256
$ws.Range("B7").Value = 0.07123092789970542
$ws.Range("D7").Value = 0.01995477856575434
$ws.Range("E7").Value = 0.1324101739341685
$ws.Range("F7").Value = 0.6264407361684903
$ws.Range("G7").Value = 0.002415947669691276
$ws.Range("I7").Value = 0.4893721152835617
$ws.Range("K7").Value = 0.4668411209017052
$ws.Range("M7").Value = 0.2249253937788112
$ws.Range("O7").Value = 2.126807704638878
$ws.Range("B8").Value = 0.08456389082384419
$ws.Range("D8").Value = 0.02212877788139167
$ws.Range("E8").Value = 0.1415582214045372
$ws.Range("F8").Value = 0.6301206345791712
$ws.Range("G8").Value = 0.002412432718377707
$ws.Range("I8").Value = 0.4787050474365149
$ws.Range("K8").Value = 0.5622696353636059
$ws.Range("M8").Value = 0.2634667971407225
$ws.Range("O8").Value = 2.121250623972827
$ws.Range("B9").Value = 0.1106327634455226
$ws.Range("D9").Value = 0.02633136336909558
$ws.Range("E9").Value = 0.1605181734519405
$ws.Range("F9").Value = 0.6411478772494519
$ws.Range("G9").Value = 0.002406232564677354
$ws.Range("I9").Value = 0.4600759883607717
$ws.Range("K9").Value = 0.747957071427237
$ws.Range("M9").Value = 0.3394197568221529
$ws.Range("O9").Value = 2.122980621899529
$ws.Range("B10").Value = 0.1297302710367489
$ws.Range("D10").Value = 0.02938097025094777
$ws.Range("E10").Value = 0.1750895851150318
$ws.Range("F10").Value = 0.6515879340774262
$ws.Range("G10").Value = 0.00240209528515251
$ws.Range("I10").Value = 0.4477863995722426
$ws.Range("K10").Value = 0.8834697191408907
$ws.Range("M10").Value = 0.3954613807714509
$ws.Range("O10").Value = 2.132000677658738
$ws.Range("B11").Value = 0.1384045801986815
$ws.Range("D11").Value = 0.03075984829850142
$ws.Range("E11").Value = 0.1818620485996689
$ws.Range("F11").Value = 0.6568484939169252
$ws.Range("G11").Value = 0.002400302975595292
$ws.Range("I11").Value = 0.4424992919968833
$ws.Range("K11").Value = 0.9449145808511332
$ws.Range("M11").Value = 0.4210108649351412
$ws.Range("O11").Value = 2.137798004057288
$ws.Range("B12").Value = 0.1416872338453317
$ws.Range("D12").Value = 0.03128076141008762
$ws.Range("E12").Value = 0.1844476015598033
$ws.Range("F12").Value = 0.6589142881366854
$ws.Range("G12").Value = 0.002399637112894718
$ws.Range("I12").Value = 0.4405408886562165
$ws.Range("K12").Value = 0.9681525774278725
$ws.Range("M12").Value = 0.4306939235938074
$ws.Range("O12").Value = 2.140237734796102
$ws.Range("B13").Value = 0.1409803532945659
$ws.Range("D13").Value = 0.03116862898564676
$ws.Range("E13").Value = 0.1838898194168763
$ws.Range("F13").Value = 0.6584660997578453
$ws.Range("G13").Value = 0.002399779948059126
$ws.Range("I13").Value = 0.4409607207706907
$ws.Range("K13").Value = 0.9631492020990322
$ws.Range("M13").Value = 0.4286081469615084
$ws.Range("O13").Value = 2.139701412158132
$ws.Range("B14").Value = 0.1386746897819364
$ws.Range("D14").Value = 0.03080272911634552
$ws.Range("E14").Value = 0.182074341751175
$ws.Range("F14").Value = 0.6570169691256211
$ws.Range("G14").Value = 0.002400247937682515
$ws.Range("I14").Value = 0.4423372965784278
$ws.Range("K14").Value = 0.9468269885616394
$ws.Range("M14").Value = 0.4218073352497385
$ws.Range("O14").Value = 2.137993819243661
$ws.Range("B15").Value = 0.1372621220041879
$ws.Range("D15").Value = 0.03057844276850119
$ws.Range("E15").Value = 0.1809650482888046
$ws.Range("F15").Value = 0.6561389432701503
$ws.Range("G15").Value = 0.00240053626584546
$ws.Range("I15").Value = 0.4431861833156567
$ws.Range("K15").Value = 0.9368252458667428
$ws.Range("M15").Value = 0.4176426863217699
$ws.Range("O15").Value = 2.136979722432329
$ws.Range("B16").Value = 0.1291630996721125
$ws.Range("D16").Value = 0.02929068568764137
$ws.Range("E16").Value = 0.1746499075119274
$ws.Range("F16").Value = 0.6512544509641884
$ws.Range("G16").Value = 0.002402214217623121
$ws.Range("I16").Value = 0.4481380372772374
$ws.Range("K16").Value = 0.8794500319102951
$ws.Range("M16").Value = 0.3937927851972916
$ws.Range("O16").Value = 2.131655966419544
$ws.Range("B17").Value = 0.1241910668260857
$ws.Range("D17").Value = 0.02849851493554212
$ws.Range("E17").Value = 0.170812818313351
$ws.Range("F17").Value = 0.6483890830760259
$ws.Range("G17").Value = 0.002403266532606303
$ws.Range("I17").Value = 0.4512536183944258
$ws.Range("K17").Value = 0.8442001213927597
$ws.Range("M17").Value = 0.3791759448198064
$ws.Range("O17").Value = 2.128824477104502
$ws.Range("B18").Value = 0.1213300513082061
$ws.Range("D18").Value = 0.02804209006218628
$ws.Range("E18").Value = 0.168619349554227
$ws.Range("F18").Value = 0.6467891153650385
$ws.Range("G18").Value = 0.00240388024863609
$ws.Range("I18").Value = 0.4530741798711126
$ws.Range("K18").Value = 0.8239064907536715
$ws.Range("M18").Value = 0.3707739839660746
$ws.Range("O18").Value = 2.127355283998412
$ws.Range("B19").Value = 0.1203611552810173
$ws.Range("D19").Value = 0.02788741780405957
$ws.Range("E19").Value = 0.167878992809726
$ws.Range("F19").Value = 0.646255651475073
$ws.Range("G19").Value = 0.002404089495800446
$ws.Range("I19").Value = 0.4536954941147187
$ws.Range("K19").Value = 0.8170322154291227
$ws.Range("M19").Value = 0.3679301279811895
$ws.Range("O19").Value = 2.126885193798131
$ws.Range("B20").Value = 0.1247204774559663
$ws.Range("D20").Value = 0.02858292475976043
$ws.Range("E20").Value = 0.1712198813432479
$ws.Range("F20").Value = 0.6486891250426368
$ws.Range("G20").Value = 0.002403153637569476
$ws.Range("I20").Value = 0.4509190029501848
$ws.Range("K20").Value = 0.8479544890654154
$ws.Range("M20").Value = 0.3807313882120411
$ws.Range("O20").Value = 2.129109390344524
$ws.Range("B21").Value = 0.1393519782580199
$ws.Range("D21").Value = 0.03091023661444581
$ws.Range("E21").Value = 0.1826070202347481
$ws.Range("F21").Value = 0.6574406111502498
$ws.Range("G21").Value = 0.002400110129600129
$ws.Range("I21").Value = 0.4419317760842079
$ws.Range("K21").Value = 0.9516220361257126
$ws.Range("M21").Value = 0.4238046808234657
$ws.Range("O21").Value = 2.138488741071143
$ws.Range("B22").Value = 0.148902060151741
$ws.Range("D22").Value = 0.03242403933091964
$ws.Range("E22").Value = 0.1901715226250289
$ws.Range("F22").Value = 0.663590088656477
$ws.Range("G22").Value = 0.002398195869796474
$ws.Range("I22").Value = 0.4363129084294517
$ws.Range("K22").Value = 1.019200385136571
$ws.Range("M22").Value = 0.4520023308373311
$ws.Range("O22").Value = 2.146043664732616
$ws.Range("B23").Value = 0.1438062114484353
$ws.Range("D23").Value = 0.0316167656452393
$ws.Range("E23").Value = 0.1861229171435212
$ws.Range("F23").Value = 0.6602685952273788
$ws.Range("G23").Value = 0.002399210717768048
$ws.Range("I23").Value = 0.43928846654988
$ws.Range("K23").Value = 0.9831488347811899
$ws.Range("M23").Value = 0.4369484385413074
$ws.Range("O23").Value = 2.141880808644743
$ws.Range("B24").Value = 0.1244811387674645
$ws.Range("D24").Value = 0.02854476617643087
$ws.Range("E24").Value = 0.171035809157317
$ws.Range("F24").Value = 0.6485533285277256
$ws.Range("G24").Value = 0.002403204650279531
$ws.Range("I24").Value = 0.4510701911039945
$ws.Range("K24").Value = 0.8462572265465838
$ws.Range("M24").Value = 0.3800281677251931
$ws.Range("O24").Value = 2.128980086880517
$ws.Range("B25").Value = 0.1035895559161162
$ws.Range("D25").Value = 0.02520105454286181
$ws.Range("E25").Value = 0.1552777409922541
$ws.Range("F25").Value = 0.6377551639616712
$ws.Range("G25").Value = 0.002407836161884142
$ws.Range("I25").Value = 0.4648704089368731
$ws.Range("K25").Value = 0.747957071427237
$ws.Range("M25").Value = 0.3188313266407121
$ws.Range("O25").Value = 2.121155482261599
